$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Elements")

# 1. Delete the "Reference.type" element row (sheet row 6), shifting subsequent rows up.
$ws.Rows(6).Delete()

# 2. Apply copy-editing fixes to cell text (typo/wording corrections).
$ws.Range("M2").Value = 'References SHALL be a reference to an actual FHIR resource, and SHALL be resolveable (allowing for access control, temporary unavailability, etc). Resolution can be either by retrieval from the URL, or, where applicable by resource type, by treating an absolute reference as a canonical URL and looking it up in a local registry/repository.'
$ws.Range("AI2").Value = 'ele-1:All FHIR elements must have a @value or children {hasValue() | (children().count() > id.count())}
ref-1:SHALL have a contained resource if a local reference is provided {reference.startsWith(''#'').not() or (reference.substring(1).trace(''url'') in %resource.contained.id.trace(''ids''))}'
$ws.Range("K3").Value = 'xml:id (or equivalent in JSON)'
$ws.Range("L3").Value = 'unique id for the element within a resource (for internal references). This may be any string value that does not contain spaces.'
$ws.Range("K4").Value = 'Additional Content defined by implementations'
$ws.Range("L4").Value = 'May be used to represent additional information that is not part of the basic definition of the element. In order to make the use of extensions safe and manageable, there is a strict set of governance  applied to the definition and use of extensions. Though any implementer is allowed to define an extension, there is a set of requirements that SHALL be met as part of the definition of the extension.'
$ws.Range("L6").Value = 'An identifier for the other resource. This is used when there is no way to reference the other resource directly, either because the entity is not available through a FHIR server, or because there is no way for the author of the resource to convert a known identifier to an actual location. There is no requirement that a Reference.identifier point to something that is actually exposed as a FHIR instance, but it SHALL point to a business concept that would be expected to be exposed as a FHIR instance, and that instance would need to be of a FHIR resource type allowed by the reference.'
$ws.Range("M6").Value = 'When an identifier is provided in place of a reference, any system processing the reference will only be able to resolve the identifier to a reference if it understands the business context in which the identifier is used. Sometimes this is global (e.g. a national identifier) but often it is not. For this reason, none of the useful mechanisms described for working with references (e.g. chaining, includes) are possible, nor should servers be expected to be able resolve the reference. Servers may accept an identifier based reference untouched, resolve it, and/or reject it - see CapabilityStatement.rest.resource.referencePolicy. 

When both an identifier and a literal reference are provided, the literal reference is preferred. Applications processing the resource are allowed - but not required - to check that the identifier matches the literal reference

Applications converting a logical reference to a literal reference may choose to leave the logical reference present, or remove it.'

# 3. Re-apply the AutoFilter over the shrunk range with the original filter criteria.
$ws.AutoFilterMode = $False
$ws.Range("A1:AJ7").AutoFilter()
$ws.Range("A1:AJ7").AutoFilter(7, "<> ")
$ws.Range("A1:AJ7").AutoFilter(27, @(""), 7)

# 4. Shrink the conditional-formatting "applies to" range to match the new row count,
#    preserving the existing rules/dxf styles instead of recreating them.
$cfRange = $ws.Range("A2:AI7")
$fcs = $cfRange.FormatConditions
for ($i = 1; $i -le $fcs.Count(); $i++) {
    $fc = $fcs.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("A2:AI6"))
}

# 5. The hidden "_FilterDatabase" defined name created by AutoFilter still points at the
#    old (pre-delete) range; repoint it at the new, smaller range.
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    if ($n.Name() -like "*_FilterDatabase*") {
        $n.RefersTo = "=Elements!`$A`$1:`$AJ`$7"
    }
}
